$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: Hydro-Dam -> Nuclear ---
$ws.Range("A5").Value = "Nuclear"
$ws.Range("B5").Value = 9787
$ws.Range("C5").Value = 135
$ws.Range("D5").Value = 2.5
$ws.Range("E5").Value = 0.93
$ws.Range("F5").Value = 10.461
$ws.Range("G5").Value = 0.66
$ws.Range("H5").Value = 0
$ws.Range("K5").Value = 0

# --- Row 6: Hydro-RoR -> Coal ---
$ws.Range("A6").Value = "Coal"
$ws.Range("B6").Formula = "=3549+693+2857"
$ws.Range("C6").Value = 78
$ws.Range("D6").Value = 8.4600000000000009
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 8.49
$ws.Range("G6").Value = 2.37
$ws.Range("H6").Value = 91.761661600000011
$ws.Range("K6").Value = 0.6

# --- Row 7: Nuclear -> Gas-CC ---
$ws.Range("A7").Value = "Gas-CC"
$ws.Range("B7").Value = 1265
$ws.Range("C7").Value = 31
$ws.Range("D7").Value = 1.96
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 6.28
$ws.Range("G7").Value = 1.6
$ws.Range("H7").Value = 53.796011199999995
$ws.Range("K7").Value = 1

# --- Row 8: Biomass -> Gas-CT ---
$ws.Range("A8").Value = "Gas-CT"
$ws.Range("B8").Value = 1120
$ws.Range("C8").Value = 24
$ws.Range("D8").Value = 6.44
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 9.7170000000000005
$ws.Range("G8").Value = 1.6
$ws.Range("H8").Value = 53.977448000000003
$ws.Range("K8").Value = 1

# --- Row 9: Coal -> Gas-CCS-95 ---
$ws.Range("A9").Value = "Gas-CCS-95"
$ws.Range("B9").Value = 2596
$ws.Range("C9").Value = 62
$ws.Range("D9").Value = 4.5199999999999996
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 7.11
$ws.Range("G9").ClearContents()
$ws.Range("H9").Value = 5.9
$ws.Range("K9").Value = 1

# --- Row 10: Gas-CC -> Gas-CCS-97 ---
$ws.Range("A10").Value = "Gas-CCS-97"
$ws.Range("B10").Value = 2635
$ws.Range("C10").Value = 62
$ws.Range("D10").Value = 4.5999999999999996
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 7.15
$ws.Range("G10").ClearContents()
$ws.Range("H10").Value = 3.6
$ws.Range("K10").Value = 1

# --- Rows 2, 3, 4: add new Ramping (K) values ---
$ws.Range("K2").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("K4").Value = 1

# --- Remove old Row 11 (Gas-CT), shifting nothing else up since row 11 is last ---
$ws.Rows("11:11").Delete()

# Re-establish the shared formulas for columns I and J across I3:I10 / J3:J10
# (row deletion above breaks the shared-formula grouping, so re-apply it)
$ws.Range("I3:I10").Formula = "=PMT(5%,30,-B3*1000)"
$ws.Range("J3:J10").Formula = "=C3*1000"

# --- New column K header ---
$ws.Range("K1").Value = "Ramping"

# Give the new Ramping column's K2 cell a plain "General" number format
# (mirrors the author clearing the inherited currency format on that cell)
$ws.Range("K2").NumberFormat = "General"

# --- Update the selected cell to match the saved view state ---
$ws.Range("I2").Select()
